# Append the next report row (row 3) to the "Item Data" sheet, mirroring
# the existing row 2 data/format, and extend the used range accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 45408
$ws.Range("E3").NumberFormat = "yyyy-MM-dd"
